$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.715.91'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.634.23'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '1.861.99'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '1.655.39'
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = '26.692.55'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '0.0₃0719'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.58%  '
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.90%  '
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("E30").Value = '  -3.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.49%  '
$ws.Range("D34").Value = '1.262.41'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("E38").Value = '  -3.83%  '
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E41").Value = '  -1.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.62%  '
$ws.Range("D43").Value = '1.772.16'
$ws.Range("E44").Value = '  -3.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  -2.11%  '
